$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The old column B ("Light to switch" data) slides over to column C, and a
# brand-new column B ("Brightness sensor" data) is introduced in its place.
# Copy column B's formatting into column C first so the relocated values
# keep their original look.
$ws.Range("B1:B9").Copy()
$ws.Range("C1").PasteSpecial(-4122)

# Fill in the new "Brightness sensor" column B top to bottom.
$ws.Range("B1").Value = "Brightness sensor"
$ws.Range("B2").Value = "A_Brightness_Sensor"
$ws.Range("B3").Value = "B_Brightness_Sensor"
$ws.Range("B4").Value = "D_Brightness_Sensor"
$ws.Range("B5").Value = "D_Brightness_Sensor"
$ws.Range("B6").Value = "E_Brightness_Sensor"
$ws.Range("B7").Value = "F_Brightness_Sensor"
$ws.Range("B8").Value = "F_Brightness_Sensor"
$ws.Range("B9").Value = "G_Brightness_Sensor"

# Column A keeps the Motion detector id (rows 5 and 8 are brand new rows
# duplicating the D and F detectors for their second analog output).
$ws.Range("A2").Value = "A_Motion_Detector"
$ws.Range("A3").Value = "B_Motion_Detector"
$ws.Range("A4").Value = "D_Motion_Detector"
$ws.Range("A5").Value = "D_Motion_Detector"
$ws.Range("A6").Value = "E_Motion_Detector"
$ws.Range("A7").Value = "F_Motion_Detector"
$ws.Range("A8").Value = "F_Motion_Detector"
$ws.Range("A9").Value = "G_Motion_Detector"

# Column C gets the analog/light id that used to live in column B, plus the
# two brand-new split rows.
$ws.Range("C1").Value = "Light to switch"
$ws.Range("C2").Value = "A_Lights_Analog"
$ws.Range("C3").Value = "B_Lights_1_Analog"
$ws.Range("C4").Value = "D_Lights_1_Analog"
$ws.Range("C5").Value = "D_Lights_2_Analog"
$ws.Range("C6").Value = "E_Lights_Analog"
$ws.Range("C7").Value = "F_Lights_1_Analog"
$ws.Range("C8").Value = "F_Lights_2_Analog"
$ws.Range("C9").Value = "G_Lights_Analog"

# Column B now holds the new (shorter) ids, so give it the plain/default
# look again, except the header cell which keeps the header style.
$ws.Range("Z1").Copy()
$ws.Range("B2:B9").PasteSpecial(-4122)
$ws.Range("Z1").ClearContents()

# Widen column B to fit the new, longer "Brightness sensor" labels
$ws.Columns.Item(2).ColumnWidth = 19

# Move the active selection to C7, matching the saved view state
$ws.Range("C7").Select()
